$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6002  # was 6501
$ws.Range("I64").Value = 6002  # was 6501
$ws.Range("K64").Value = 6002  # was 6501
$ws.Range("M64").Value = -5754  # was -6253
$ws.Range("H67").Value = 6002  # was 6501
$ws.Range("I67").Value = 6002  # was 6501
$ws.Range("K67").Value = 6002  # was 6501
$ws.Range("M67").Value = -5144  # was -5643
$ws.Range("H98").Value = 1863.1428  # was 1594.2
$ws.Range("I98").Value = 1863.1428  # was 1594.2
$ws.Range("K98").Value = 1863.1428  # was 1594.2
$ws.Range("M98").Value = -365.1428000000001  # was -96.20000000000005
$ws.Range("H100").Value = 1838.8125  # was 2612.625
$ws.Range("I100").Value = 1324.8462  # was 1740.6
$ws.Range("K100").Value = 1324.8462  # was 1740.6
$ws.Range("M100").Value = -783.8462  # was -1199.6
$ws.Range("H116").Value = 5846  # was 6369.8096
$ws.Range("I116").Value = 5379.2856  # was 6377.2
$ws.Range("J116").Value = 6499.4  # was 6363.091
$ws.Range("K116").Value = 5379.2856  # was 6377.2
$ws.Range("L116").Value = 6499.4  # was 6363.091
$ws.Range("M116").Value = -1937.2856  # was -2935.2
$ws.Range("N116").Value = -13383.4  # was -13247.091
$ws.Range("H122").Value = 1863.1428  # was 1594.2
$ws.Range("I122").Value = 1863.1428  # was 1594.2
$ws.Range("K122").Value = 5589.428400000001  # was 4782.6
$ws.Range("M122").Value = -3139.428400000001  # was -2332.6
$ws.Range("H132").Value = 4014.2163  # was 3902.5676
$ws.Range("I132").Value = 3192.3235  # was 3125.8857
$ws.Range("J132").Value = 13329  # was 17494.5
$ws.Range("K132").Value = 9576.970499999999  # was 9377.6571
$ws.Range("L132").Value = 39987  # was 52483.5
$ws.Range("M132").Value = -7046.970499999999  # was -6847.6571
$ws.Range("N132").Value = -45047  # was -57543.5
$ws.Range("H135").Value = 1775.6  # was 1928.826
$ws.Range("I135").Value = 1244.7778  # was 1339.1177
$ws.Range("J135").Value = 3140.5715  # was 3599.6667
$ws.Range("K135").Value = 11203.0002  # was 12052.0593
$ws.Range("L135").Value = 28265.1435  # was 32397.0003
$ws.Range("M135").Value = -8668.0002  # was -9517.059300000001
$ws.Range("N135").Value = -33335.1435  # was -37467.0003

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1829.2354  # was 1794.0555
$ws.Range("I61").Value = 1693.5625  # was 1664.2941
$ws.Range("K61").Value = 1693.5625  # was 1664.2941
$ws.Range("M61").Value = -1481.5625  # was -1452.2941
$ws.Range("H74").Value = 53992.105  # was 60208.65
$ws.Range("I74").Value = 63739.812  # was 67935.53
$ws.Range("J74").Value = 2004.3334  # was 2257
$ws.Range("K74").Value = 63739.812  # was 67935.53
$ws.Range("L74").Value = 2004.3334  # was 2257
$ws.Range("M74").Value = -62865.812  # was -67061.53
$ws.Range("N74").Value = -3752.3334  # was -4005
$ws.Range("H77").Value = 53992.105  # was 60208.65
$ws.Range("I77").Value = 63739.812  # was 67935.53
$ws.Range("J77").Value = 2004.3334  # was 2257
$ws.Range("K77").Value = 318699.06  # was 339677.65
$ws.Range("L77").Value = 10021.667  # was 11285
$ws.Range("M77").Value = -314331.06  # was -335309.65
$ws.Range("N77").Value = -18757.667  # was -20021
$ws.Range("H80").Value = 100000  # was 77499.25
$ws.Range("I80").Value = 0  # was 79999
$ws.Range("J80").Value = 100000  # was 76666
$ws.Range("K80").Value = 0  # was 79999
$ws.Range("L80").Value = 100000  # was 76666
$ws.Range("M80").ClearContents()  # was -79001
$ws.Range("N80").Value = -101996  # was -78662
$ws.Range("H83").Value = 100000  # was 77499.25
$ws.Range("I83").Value = 0  # was 79999
$ws.Range("J83").Value = 100000  # was 76666
$ws.Range("K83").Value = 0  # was 239997
$ws.Range("L83").Value = 300000  # was 229998
$ws.Range("M83").ClearContents()  # was -235005
$ws.Range("N83").Value = -309984  # was -239982
$ws.Range("H136").Value = 1829.2354  # was 1794.0555
$ws.Range("I136").Value = 1693.5625  # was 1664.2941
$ws.Range("K136").Value = 5080.6875  # was 4992.8823
$ws.Range("M136").Value = -2530.6875  # was -2442.8823

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 3321.3333  # was 3699.375
$ws.Range("I64").Value = 1719.6  # was 1739.6
$ws.Range("J64").Value = 5323.5  # was 6965.6665
$ws.Range("K64").Value = 1719.6  # was 1739.6
$ws.Range("L64").Value = 5323.5  # was 6965.6665
$ws.Range("M64").Value = -1494.6  # was -1514.6
$ws.Range("N64").Value = -5773.5  # was -7415.6665
$ws.Range("H67").Value = 3321.3333  # was 3699.375
$ws.Range("I67").Value = 1719.6  # was 1739.6
$ws.Range("J67").Value = 5323.5  # was 6965.6665
$ws.Range("K67").Value = 1719.6  # was 1739.6
$ws.Range("L67").Value = 5323.5  # was 6965.6665
$ws.Range("M67").Value = -939.5999999999999  # was -959.5999999999999
$ws.Range("N67").Value = -6883.5  # was -8525.666499999999
$ws.Range("H86").Value = 2996.3333  # was 2597.6
$ws.Range("I86").Value = 0  # was 1500
$ws.Range("J86").Value = 2996.3333  # was 2872
$ws.Range("K86").Value = 0  # was 1500
$ws.Range("L86").Value = 2996.3333  # was 2872
$ws.Range("M86").ClearContents()  # was -377
$ws.Range("N86").Value = -5242.3333  # was -5118
$ws.Range("H89").Value = 2996.3333  # was 2597.6
$ws.Range("I89").Value = 0  # was 1500
$ws.Range("J89").Value = 2996.3333  # was 2872
$ws.Range("K89").Value = 0  # was 7500
$ws.Range("L89").Value = 14981.6665  # was 14360
$ws.Range("M89").ClearContents()  # was -1884
$ws.Range("N89").Value = -26213.6665  # was -25592
$ws.Range("H105").Value = 4077.3667  # was 4328.875
$ws.Range("I105").Value = 3805.375  # was 4135.846
$ws.Range("K105").Value = 3805.375  # was 4135.846
$ws.Range("M105").Value = -2058.375  # was -2388.846
$ws.Range("H134").Value = 1794.1428  # was 1762.9546
$ws.Range("I134").Value = 1736.05  # was 1706.1428
$ws.Range("K134").Value = 5208.15  # was 5118.428400000001
$ws.Range("M134").Value = -2673.15  # was -2583.428400000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 839.9375  # was 757.4761999999999
$ws.Range("I16").Value = 795.93335  # was 720.35
$ws.Range("K16").Value = 795.93335  # was 720.35
$ws.Range("M16").Value = -508.93335  # was -433.35
$ws.Range("H31").Value = 4122.759  # was 4498.577
$ws.Range("I31").Value = 1994  # was 2219.6667
$ws.Range("K31").Value = 1994  # was 2219.6667
$ws.Range("M31").Value = -1699  # was -1924.6667
$ws.Range("H34").Value = 4122.759  # was 4498.577
$ws.Range("I34").Value = 1994  # was 2219.6667
$ws.Range("K34").Value = 1994  # was 2219.6667
$ws.Range("M34").Value = -1792  # was -2017.6667
$ws.Range("H58").Value = 52892.35  # was 50608.24
$ws.Range("I58").Value = 58468.555  # was 55650.527
$ws.Range("K58").Value = 58468.555  # was 55650.527
$ws.Range("M58").Value = -58265.555  # was -55447.527
$ws.Range("H62").Value = 3900  # was 3266.6667
$ws.Range("J62").Value = 4000  # was 3000
$ws.Range("L62").Value = 4000  # was 3000
$ws.Range("N62").Value = -5248  # was -4248
$ws.Range("H65").Value = 3900  # was 3266.6667
$ws.Range("J65").Value = 4000  # was 3000
$ws.Range("L65").Value = 20000  # was 15000
$ws.Range("N65").Value = -26240  # was -21240
$ws.Range("H68").Value = 51733  # was 27000
$ws.Range("J68").Value = 64099.5  # was 0
$ws.Range("L68").Value = 64099.5  # was 0
$ws.Range("N68").Value = -65597.5  # was None
$ws.Range("H71").Value = 51733  # was 27000
$ws.Range("J71").Value = 64099.5  # was 0
$ws.Range("L71").Value = 192298.5  # was 0
$ws.Range("N71").Value = -199786.5  # was None
$ws.Range("H74").Value = 0  # was 30000
$ws.Range("J74").Value = 0  # was 30000
$ws.Range("L74").Value = 0  # was 30000
$ws.Range("N74").ClearContents()  # was -31748
$ws.Range("H77").Value = 0  # was 30000
$ws.Range("J77").Value = 0  # was 30000
$ws.Range("L77").Value = 0  # was 90000
$ws.Range("N77").ClearContents()  # was -98736
$ws.Range("H113").Value = 839.9375  # was 757.4761999999999
$ws.Range("I113").Value = 795.93335  # was 720.35
$ws.Range("K113").Value = 795.93335  # was 720.35
$ws.Range("M113").Value = 1374.06665  # was 1449.65
$ws.Range("H136").Value = 52892.35  # was 50608.24
$ws.Range("I136").Value = 58468.555  # was 55650.527
$ws.Range("K136").Value = 175405.665  # was 166951.581
$ws.Range("M136").Value = -172855.665  # was -164401.581

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4772036.5  # was 4357293.5
$ws.Range("J131").Value = 8349095  # was 7156723.5
$ws.Range("L131").Value = 25047285  # was 21470170.5
$ws.Range("N131").Value = -25057365  # was -21480250.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 925.9  # was 994
$ws.Range("I97").Value = 977.1875  # was 1040.3572
$ws.Range("J97").Value = 720.75  # was 777.6667
$ws.Range("K97").Value = 977.1875  # was 1040.3572
$ws.Range("L97").Value = 720.75  # was 777.6667
$ws.Range("M97").Value = -481.1875  # was -544.3571999999999
$ws.Range("N97").Value = -1712.75  # was -1769.6667
$ws.Range("H107").Value = 67862.8  # was 72613.64
$ws.Range("I107").Value = 125443  # was 143349.14
$ws.Range("J107").Value = 2056.8572  # was 1878.1428
$ws.Range("K107").Value = 125443  # was 143349.14
$ws.Range("L107").Value = 2056.8572  # was 1878.1428
$ws.Range("M107").Value = -123523  # was -141429.14
$ws.Range("N107").Value = -5896.8572  # was -5718.1428
$ws.Range("H122").Value = 3497.5  # was 3496.6667
$ws.Range("I122").Value = 3497.5  # was 3496.6667
$ws.Range("K122").Value = 10492.5  # was 10490.0001
$ws.Range("M122").Value = -8042.5  # was -8040.000100000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5016.75  # was 5997.5884
$ws.Range("I7").Value = 5338.25  # was 5671.857
$ws.Range("J7").Value = 4695.25  # was 6225.6
$ws.Range("K7").Value = 5338.25  # was 5671.857
$ws.Range("L7").Value = 4695.25  # was 6225.6
$ws.Range("M7").Value = -5226.25  # was -5559.857
$ws.Range("N7").Value = -4919.25  # was -6449.6
$ws.Range("H40").Value = 4354.3335  # was 4355
$ws.Range("I40").Value = 4354  # was 0
$ws.Range("K40").Value = 4354  # was 0
$ws.Range("M40").Value = -4218  # was None
$ws.Range("H68").Value = 7000  # was 6999.5
$ws.Range("I68").Value = 0  # was 6999.5
$ws.Range("J68").Value = 7000  # was 0
$ws.Range("K68").Value = 0  # was 6999.5
$ws.Range("L68").Value = 7000  # was 0
$ws.Range("M68").ClearContents()  # was -6250.5
$ws.Range("N68").Value = -8498  # was None
$ws.Range("H71").Value = 7000  # was 6999.5
$ws.Range("I71").Value = 0  # was 6999.5
$ws.Range("J71").Value = 7000  # was 0
$ws.Range("K71").Value = 0  # was 34997.5
$ws.Range("L71").Value = 35000  # was 0
$ws.Range("M71").ClearContents()  # was -31253.5
$ws.Range("N71").Value = -42488  # was None
$ws.Range("H82").Value = 2334.875  # was 2524.9656
$ws.Range("J82").Value = 2710.261  # was 3042.2
$ws.Range("L82").Value = 2710.261  # was 3042.2
$ws.Range("N82").Value = -3432.261  # was -3764.2
$ws.Range("H85").Value = 2334.875  # was 2524.9656
$ws.Range("J85").Value = 2710.261  # was 3042.2
$ws.Range("L85").Value = 2710.261  # was 3042.2
$ws.Range("N85").Value = -5206.261  # was -5538.2
$ws.Range("H122").Value = 4066.3948  # was 3997.0854
$ws.Range("I122").Value = 3423.4375  # was 3379.054
$ws.Range("J122").Value = 4534  # was 4505.2446
$ws.Range("K122").Value = 10270.3125  # was 10137.162
$ws.Range("L122").Value = 13602  # was 13515.7338
$ws.Range("M122").Value = -7820.3125  # was -7687.162
$ws.Range("N122").Value = -18502  # was -18415.7338
$ws.Range("H126").Value = 5016.75  # was 5997.5884
$ws.Range("I126").Value = 5338.25  # was 5671.857
$ws.Range("J126").Value = 4695.25  # was 6225.6
$ws.Range("K126").Value = 16014.75  # was 17015.571
$ws.Range("L126").Value = 14085.75  # was 18676.8
$ws.Range("M126").Value = -13544.75  # was -14545.571
$ws.Range("N126").Value = -19025.75  # was -23616.8
$ws.Range("H132").Value = 40075.438  # was 37878.293
$ws.Range("I132").Value = 45118.715  # was 43602.516
$ws.Range("J132").Value = 4772.5  # was 4677.8
$ws.Range("K132").Value = 135356.145  # was 130807.548
$ws.Range("L132").Value = 14317.5  # was 14033.4
$ws.Range("M132").Value = -132826.145  # was -128277.548
$ws.Range("N132").Value = -19377.5  # was -19093.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 68681.74000000001  # was 68818.85000000001
$ws.Range("I132").Value = 72187.14  # was 72337.31
$ws.Range("K132").Value = 216561.42  # was 217011.93
$ws.Range("M132").Value = -214031.42  # was -214481.93
$ws.Range("H133").Value = 89999  # was 0
$ws.Range("J133").Value = 89999  # was 0
$ws.Range("L133").Value = 89999  # was 0
$ws.Range("N133").Value = -100119  # was None
